$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Homeworks")

# Hm06 (column G) correction for student row 23
$ws.Cells.Item(23, 7).Value = 34

# Newly recorded Hm07 (column H) and Hm08 (column I) scores
$hiData = @{
  2  = @(0, 0)
  3  = @(56, 68)
  4  = @(60, 70)
  5  = @(0, 0)
  6  = @(0, 0)
  7  = @(57, 57.5)
  8  = @(0, 3)
  9  = @(17, 0)
  10 = @(58.5, 67)
  11 = @(0, 0)
  12 = @(0, 66)
  13 = @(52, 59)
  14 = @(30, 0)
  15 = @(38.5, 54)
  16 = @(48, 48)
  17 = @(50.5, 56.5)
  18 = @(0, 0)
  19 = @(53.5, 59)
  20 = @(54, 64)
  21 = @(30, 27)
  22 = @(60, 67)
  23 = @(42.5, 0)
  24 = @(22, 48.5)
  25 = @(31, 43)
  26 = @(44.5, 68)
  27 = @(24, 44)
  28 = @(0, 17.5)
  29 = @(18, 31.5)
  30 = @(53.5, 62)
  31 = @(39.5, 48)
  32 = @(42, 44.5)
  33 = @(0, 0)
  34 = @(28, 44.5)
  35 = @(0, 9)
  36 = @(0, 0)
  37 = @(0, 0)
  38 = @(35, 58)
  39 = @(0, 33)
  40 = @(31, 45)
  41 = @(60, 64)
  42 = @(42, 56.5)
  43 = @(59.5, 59.5)
  44 = @(52.5, 67.5)
  45 = @(0, 0)
  46 = @(50, 45)
  47 = @(0, 0)
  48 = @(0, 0)
}

foreach ($r in 2..48) {
  $vals = $hiData[$r]
  $hCell = $ws.Cells.Item($r, 8)
  $iCell = $ws.Cells.Item($r, 9)
  $hCell.Value = $vals[0]
  $hCell.NumberFormat = "0.00"
  $iCell.Value = $vals[1]
  $iCell.NumberFormat = "0.00"
}

$excel.CalculateFull()
